$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Label" header in H1, copying the header style/format from G1
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Label"

# Populate the new Label column (0 = not the aggregate/summary row, 1 = summary row)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("H7").Value = 1
$ws.Range("H8").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("H10").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("H12").Value = 0
$ws.Range("H13").Value = 1

# Update refitted prediction / error / cross-entropy values
$ws.Range("D2").Value = 0.6052147021520106
$ws.Range("E2").Value = 0.6052147021520106

$ws.Range("D3").Value = 0.5678687588387376
$ws.Range("E3").Value = 0.5678687588387376

$ws.Range("D4").Value = 0.4489620131221874
$ws.Range("E4").Value = 0.4489620131221874

$ws.Range("D5").Value = 0.2400759035473312
$ws.Range("E5").Value = 0.2400759035473312

$ws.Range("D6").Value = 0.5922896008012222
$ws.Range("E6").Value = 0.5922896008012222

$ws.Range("F7").Value = 0.6591851711273193
